$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.887.33"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "1.709.97"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.54%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "318.71"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("E6").Value = "  +0.50%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3971"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.51%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4109"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +2.17%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "1.514"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -0.74%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +0.69%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "53.32"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.42%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.09001"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +2.61%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.712"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +6.57%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "24.38"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +4.58%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.00001395"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +6.14%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "8.219"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("D17").Value = "1.727.47"
$ws.Range("E17").Value = "  +1.73%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "100.28"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.47%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.07175"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.47%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.585"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +7.51%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "20.12"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("E22").Value = "  +0.86%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "14.56"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D24").Value = "24.871.77"
$ws.Range("E24").Value = "  +1.04%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.094"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.61%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.346"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.48%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "23.07"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +0.96%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.303"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +23.92%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "166.54"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +2.79%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "139.53"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +2.15%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "5.240"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +0.75%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.915"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +10.48%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09156"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +6.42%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.087"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +0.54%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.03061"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +11.85%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.2822"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("E37").Value = "  -3.27%  "

$ws.Range("E38").Value = "  +2.47%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "14.63"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +1.26%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.09344"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.35%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.490"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +0.06%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.7882"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +2.96%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "16.72"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +7.18%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.653"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +2.09%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.7345"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +2.30%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "4.269"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("B47").Value = "Flow"
$ws.Range("C47").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.359"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +2.92%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.41%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "141.31"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.14%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "94.88"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +4.96%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.08081"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +1.15%  "

